$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 1114 (shifts the existing rows 1114:1198 down to 1117:1201,
# and extends the used range to A1:T1201, matching the diff's new <dimension>).
$ws.Rows("1114:1116").Insert()

# Fill the 3 newly inserted rows with the new weekly price entries (date 44714 = 2022-06-02)
$newRows = New-Object 'object[,]' 3,20
# Row 1 -> sheet row 1114
$newRows[0,0]  = 9
$newRows[0,1]  = "Vega Central Mapocho de Santiago"
$newRows[0,2]  = "Metropolitana"
$newRows[0,3]  = 44714
$newRows[0,4]  = 13
$newRows[0,5]  = "Fruta"
$newRows[0,6]  = 100108
$newRows[0,7]  = "Tropicales y subtropicales"
$newRows[0,8]  = 100108006
$newRows[0,9]  = "Plátano"
$newRows[0,10] = "Sin especificar"
$newRows[0,11] = "Pintón"
$newRows[0,12] = 910
$newRows[0,13] = 9000
$newRows[0,14] = 9000
$newRows[0,15] = 9000
$newRows[0,16] = "`$/caja 20 kilos"
$newRows[0,17] = "Ecuador"
$newRows[0,18] = 450
$newRows[0,19] = 20

# Row 2 -> sheet row 1115
$newRows[1,0]  = 9
$newRows[1,1]  = "Vega Central Mapocho de Santiago"
$newRows[1,2]  = "Metropolitana"
$newRows[1,3]  = 44714
$newRows[1,4]  = 13
$newRows[1,5]  = "Fruta"
$newRows[1,6]  = 100108
$newRows[1,7]  = "Tropicales y subtropicales"
$newRows[1,8]  = 100108006
$newRows[1,9]  = "Plátano"
$newRows[1,10] = "Sin especificar"
$newRows[1,11] = "Primera Maduro"
$newRows[1,12] = 780
$newRows[1,13] = 11000
$newRows[1,14] = 11000
$newRows[1,15] = 11000
$newRows[1,16] = "`$/caja 20 kilos"
$newRows[1,17] = "Ecuador"
$newRows[1,18] = 550
$newRows[1,19] = 20

# Row 3 -> sheet row 1116
$newRows[2,0]  = 9
$newRows[2,1]  = "Vega Central Mapocho de Santiago"
$newRows[2,2]  = "Metropolitana"
$newRows[2,3]  = 44714
$newRows[2,4]  = 13
$newRows[2,5]  = "Fruta"
$newRows[2,6]  = 100108
$newRows[2,7]  = "Tropicales y subtropicales"
$newRows[2,8]  = 100108006
$newRows[2,9]  = "Plátano"
$newRows[2,10] = "Sin especificar"
$newRows[2,11] = "Primera Pintón"
$newRows[2,12] = 850
$newRows[2,13] = 12000
$newRows[2,14] = 12000
$newRows[2,15] = 12000
$newRows[2,16] = "`$/caja 20 kilos"
$newRows[2,17] = "Ecuador"
$newRows[2,18] = 600
$newRows[2,19] = 20

$ws.Range("A1114:T1116").Value = $newRows

Write-Output "done"
